$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-12: update B, C, D columns
$ws.Range("B2").Value = 79488
$ws.Range("C2").Value = 5597.0385765805
$ws.Range("D2").Value = 73890.9614234195

$ws.Range("B3").Value = 75184
$ws.Range("C3").Value = 5427.810091879
$ws.Range("D3").Value = 69756.189908121

$ws.Range("B4").Value = 71544
$ws.Range("C4").Value = 5378.673491737501
$ws.Range("D4").Value = 66165.3265082625

$ws.Range("B5").Value = 69610
$ws.Range("C5").Value = 5348.731399226001
$ws.Range("D5").Value = 64261.268600774

$ws.Range("B6").Value = 70806
$ws.Range("C6").Value = 5372.321170922
$ws.Range("D6").Value = 65433.678829078

$ws.Range("B7").Value = 73348
$ws.Range("C7").Value = 5496.759895401
$ws.Range("D7").Value = 67851.240104599

$ws.Range("B8").Value = 78137
$ws.Range("C8").Value = 6441.069913339
$ws.Range("D8").Value = 71695.930086661

$ws.Range("B9").Value = 80939
$ws.Range("C9").Value = 7184.668538918501
$ws.Range("D9").Value = 73754.3314610815

$ws.Range("B10").Value = 103039
$ws.Range("C10").Value = 8465.353500000001
$ws.Range("D10").Value = 94573.6465

$ws.Range("B11").Value = 114408
$ws.Range("C11").Value = 14333.5385
$ws.Range("D11").Value = 100074.4615

$ws.Range("B12").Value = 96449.10000000001
$ws.Range("C12").Value = 15915.495
$ws.Range("D12").Value = 80533.60500000001

# Rows 13-20: B column cell is removed entirely, C and D updated
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = 15771.6165
$ws.Range("D13").Value = 103987.3835

$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = 15824.06
$ws.Range("D14").Value = 76684.94

$ws.Range("B15").ClearContents()
$ws.Range("C15").Value = 15963.6205
$ws.Range("D15").Value = 103417.3795

$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = 16132.936
$ws.Range("D16").Value = 110937.064

$ws.Range("B17").ClearContents()
$ws.Range("C17").Value = 16102.835
$ws.Range("D17").Value = 87260.16500000001

$ws.Range("B18").ClearContents()
$ws.Range("C18").Value = 16787.6375
$ws.Range("D18").Value = 81550.3625

$ws.Range("B19").ClearContents()
$ws.Range("C19").Value = 16403.816
$ws.Range("D19").Value = 75216.18400000001

$ws.Range("B20").ClearContents()
$ws.Range("C20").Value = 16006.496
$ws.Range("D20").Value = 71081.504

# Rows 21-25: only C and D columns (B already empty)
$ws.Range("C21").Value = 14046.071
$ws.Range("D21").Value = 69127.929

$ws.Range("C22").Value = 11183.005
$ws.Range("D22").Value = 68523.995

$ws.Range("C23").Value = 9240.505000000001
$ws.Range("D23").Value = 66852.495

$ws.Range("C24").Value = 6857.7675
$ws.Range("D24").Value = 33754.2325

$ws.Range("C25").Value = 6096.6435
